# Updated cryptos list with GitHub Actions: refresh Price (D) and
# Volume(1h) (E) columns for the crypto rows. Price cells are forced to
# text ("@" number format, restored to the "Normal" style afterwards) so
# values such as "384.14" are not auto-converted to numeric cells by
# Excel - this matches the original inline-string/text representation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.315.24"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.034.41"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.44%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "384.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.51"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("E7").Value = "  -0.35%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.581"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.75"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.65%  "
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("E12").Value = "  +0.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.519.65"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.61"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.80%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.71"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.79%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.036.75"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.965"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.51"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "51.413.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("E20").Value = "  -0.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0964"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.55%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.12"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "267.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.14"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.18"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.80"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.26%  "
$ws.Range("E28").Value = "  +2.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.94%  "
$ws.Range("E31").Value = "  -2.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.77"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.39%  "
$ws.Range("E34").Value = "  -0.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "50.19"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.94%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0444"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.75%  "
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.32"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.28%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.287"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.94"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.89%  "
$ws.Range("E41").Value = "  +1.21%  "
$ws.Range("E42").Value = "  -0.84%  "
$ws.Range("E43").Value = "  -0.19%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "124.35"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("E45").Value = "  +4.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "21.69"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.55%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.46"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.77%  "
$ws.Range("E48").Value = "  +2.34%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.023.01"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.339.24"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.56%  "
$ws.Range("E51").Value = "  +5.54%  "
